# Update the "About" sheet's last-updated date (C1) and the "MCF" sheet's
# maximum-capacity-factor assumptions (column B), per the 4.0-files update.

$wb = $excel.ActiveWorkbook

# --- "About" sheet: bump the last-updated date from 1/29/2024 to 4/10/2024 ---
$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Range("C1").Value = 45392

# --- "MCF" sheet: raise several plant types' capacity factor to 1 (100%) ---
$wsMcf = $wb.Worksheets.Item("MCF")

$wsMcf.Range("B2").Value  = 1   # hard coal
$wsMcf.Range("B3").Value  = 1   # natural gas steam turbine
$wsMcf.Range("B4").Value  = 1   # natural gas combined cycle
$wsMcf.Range("B6").Value  = 1   # hydro
$wsMcf.Range("B10").Value = 1   # biomass
$wsMcf.Range("B11").Value = 1   # geothermal
$wsMcf.Range("B12").Value = 1   # petroleum
$wsMcf.Range("B13").Value = 1   # natural gas peaker
$wsMcf.Range("B14").Value = 1   # lignite
$wsMcf.Range("B16").Value = 1   # crude oil
$wsMcf.Range("B17").Value = 1   # heavy or residual fuel oil
$wsMcf.Range("B18").Value = 1   # municipal solid waste

# B19, B20, B21, B22, B24, B25 are formulas referencing the cells above
# (e.g. =B2, =B4, =B10, =B14) and will recalculate automatically.

# Update the active selection on the MCF sheet to match the saved view.
$wsMcf.Activate()
$wsMcf.Range("B17").Select()
